$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Nb nouveaux cas positifs" (column C) for existing rows - cumulative
# column B recalculates automatically from these.
$ws.Cells.Item(568, 3).Value = 88
$ws.Cells.Item(570, 3).Value = 56
$ws.Cells.Item(571, 3).Value = 20
$ws.Cells.Item(572, 3).Value = 25
$ws.Cells.Item(573, 3).Value = 68

# Fill in newly reported data for rows 574-576 (previously blank days)
# Columns: C=new cases, E=SI total, F=intubated, G=hosp hors SI, L=new deaths hopital, M=new deaths extra-hosp
# L/M are formatted as Text ("@"); write through a General round-trip so the
# values land as real numbers (matching the rest of the column) instead of
# being coerced to text by the cell's text number format.
$ws.Cells.Item(574, 3).Value = 58
$ws.Cells.Item(574, 5).Value = 6
$ws.Cells.Item(574, 6).Value = 2
$ws.Cells.Item(574, 7).Value = 15
$fmtL = $ws.Cells.Item(574, 12).NumberFormat
$ws.Cells.Item(574, 12).NumberFormat = "general"
$ws.Cells.Item(574, 12).Value = 0
$ws.Cells.Item(574, 12).NumberFormat = $fmtL
$fmtM = $ws.Cells.Item(574, 13).NumberFormat
$ws.Cells.Item(574, 13).NumberFormat = "general"
$ws.Cells.Item(574, 13).Value = 0
$ws.Cells.Item(574, 13).NumberFormat = $fmtM

$ws.Cells.Item(575, 3).Value = 43
$ws.Cells.Item(575, 5).Value = 6
$ws.Cells.Item(575, 6).Value = 2
$ws.Cells.Item(575, 7).Value = 13
$fmtL = $ws.Cells.Item(575, 12).NumberFormat
$ws.Cells.Item(575, 12).NumberFormat = "general"
$ws.Cells.Item(575, 12).Value = 0
$ws.Cells.Item(575, 12).NumberFormat = $fmtL
$fmtM = $ws.Cells.Item(575, 13).NumberFormat
$ws.Cells.Item(575, 13).NumberFormat = "general"
$ws.Cells.Item(575, 13).Value = 0
$ws.Cells.Item(575, 13).NumberFormat = $fmtM

$ws.Cells.Item(576, 3).Value = 3
$ws.Cells.Item(576, 5).Value = 5
$ws.Cells.Item(576, 6).Value = 2
$ws.Cells.Item(576, 7).Value = 13
$fmtL = $ws.Cells.Item(576, 12).NumberFormat
$ws.Cells.Item(576, 12).NumberFormat = "general"
$ws.Cells.Item(576, 12).Value = 0
$ws.Cells.Item(576, 12).NumberFormat = $fmtL
$fmtM = $ws.Cells.Item(576, 13).NumberFormat
$ws.Cells.Item(576, 13).NumberFormat = "general"
$ws.Cells.Item(576, 13).Value = 0
$ws.Cells.Item(576, 13).NumberFormat = $fmtM

# K573/K574 re-enter the same "new deaths today" formula explicitly (the
# author's edit broke these two out of the shared formula run si=46 while
# K575/K576 kept using it); value is unchanged (0) either way.
$ws.Cells.Item(573, 11).Formula = "=IF(TODAY()>A572,L573+M573,"""")"
$ws.Cells.Item(574, 11).Formula = "=IF(TODAY()>A573,L574+M574,"""")"

# Move the active selection on the frozen (bottom-right) pane to A2
$ws.Cells.Item(2, 1).Select()
